$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2198.889
$ws.Range("J62").Value = 2497.5
$ws.Range("L62").Value = 2497.5
$ws.Range("N62").Value = -3745.5

$ws.Range("H65").Value = 2198.889
$ws.Range("J65").Value = 2497.5
$ws.Range("L65").Value = 12487.5
$ws.Range("N65").Value = -18727.5

$ws.Range("H70").Value = 1632.3334
$ws.Range("I70").Value = 2147
$ws.Range("J70").Value = 1375
$ws.Range("K70").Value = 6441
$ws.Range("L70").Value = 4125
$ws.Range("M70").Value = -6171
$ws.Range("N70").Value = -4665

$ws.Range("H73").Value = 1632.3334
$ws.Range("I73").Value = 2147
$ws.Range("J73").Value = 1375
$ws.Range("K73").Value = 6441
$ws.Range("L73").Value = 4125
$ws.Range("M73").Value = -5505
$ws.Range("N73").Value = -5997

$ws.Range("H113").Value = 92745.91
$ws.Range("J113").Value = 1914.2858
$ws.Range("L113").Value = 1914.2858
$ws.Range("N113").Value = -8422.2858

$ws.Range("H125").Value = 2689.8928
$ws.Range("I125").Value = 2223.9443
$ws.Range("J125").Value = 3528.6
$ws.Range("K125").Value = 20015.4987
$ws.Range("L125").Value = 31757.4
$ws.Range("M125").Value = -17555.4987
$ws.Range("N125").Value = -36677.39999999999

$ws.Range("H132").Value = 3971992.8
$ws.Range("I132").Value = 4468173
$ws.Range("K132").Value = 13404519
$ws.Range("M132").Value = -13401989

$ws.Range("H135").Value = 786.44446
$ws.Range("I135").Value = 650.4666999999999
$ws.Range("J135").Value = 1466.3334
$ws.Range("K135").Value = 5854.2003
$ws.Range("L135").Value = 13197.0006
$ws.Range("M135").Value = -3319.2003
$ws.Range("N135").Value = -18267.0006

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 30633.918
$ws.Range("I32").Value = 6169.6733
$ws.Range("J32").Value = 130529.586
$ws.Range("K32").Value = 6169.6733
$ws.Range("L32").Value = 130529.586
$ws.Range("M32").Value = -5882.6733
$ws.Range("N32").Value = -131103.586

$ws.Range("H61").Value = 1950.0857
$ws.Range("I61").Value = 1905.1482
$ws.Range("J61").Value = 2101.75
$ws.Range("K61").Value = 1905.1482
$ws.Range("L61").Value = 2101.75
$ws.Range("M61").Value = -1693.1482
$ws.Range("N61").Value = -2525.75

$ws.Range("H76").Value = 30429.334
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 30429.334
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 30429.334
$ws.Range("M76").ClearContents()
$ws.Range("N76").Value = -31105.334

$ws.Range("H79").Value = 30429.334
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 30429.334
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 30429.334
$ws.Range("M79").ClearContents()
$ws.Range("N79").Value = -32769.334

$ws.Range("H136").Value = 1950.0857
$ws.Range("I136").Value = 1905.1482
$ws.Range("J136").Value = 2101.75
$ws.Range("K136").Value = 5715.444600000001
$ws.Range("L136").Value = 6305.25
$ws.Range("M136").Value = -3165.444600000001
$ws.Range("N136").Value = -11405.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H98").Value = 25000
$ws.Range("I98").Value = 25000
$ws.Range("K98").Value = 25000
$ws.Range("M98").Value = -22754

$ws.Range("H99").Value = 22458.8
$ws.Range("I99").Value = 3493.3333
$ws.Range("K99").Value = 3493.3333
$ws.Range("M99").Value = -1995.3333

$ws.Range("H126").Value = 22458.8
$ws.Range("I126").Value = 3493.3333
$ws.Range("K126").Value = 10479.9999
$ws.Range("M126").Value = -8009.999899999999

$ws.Range("H134").Value = 1122.7693
$ws.Range("I134").Value = 1052.3243
$ws.Range("J134").Value = 1296.5333
$ws.Range("K134").Value = 3156.9729
$ws.Range("L134").Value = 3889.5999
$ws.Range("M134").Value = -621.9728999999998
$ws.Range("N134").Value = -8959.599900000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 3055.6667
$ws.Range("I68").Value = 1250.5
$ws.Range("J68").Value = 6666
$ws.Range("K68").Value = 3751.5
$ws.Range("L68").Value = 19998
$ws.Range("M68").Value = -2940.5
$ws.Range("N68").Value = -21620

$ws.Range("H71").Value = 3055.6667
$ws.Range("I71").Value = 1250.5
$ws.Range("J71").Value = 6666
$ws.Range("K71").Value = 11254.5
$ws.Range("L71").Value = 59994
$ws.Range("M71").Value = -7198.5
$ws.Range("N71").Value = -68106

$ws.Range("H127").Value = 683
$ws.Range("J127").Value = 683
$ws.Range("L127").Value = 2049
$ws.Range("N127").Value = -11969

$ws.Range("H129").Value = 10206487
$ws.Range("I129").Value = 25005966
$ws.Range("J129").Value = 340167.6
$ws.Range("K129").Value = 75017898
$ws.Range("L129").Value = 1020502.8
$ws.Range("M129").Value = -75012898
$ws.Range("N129").Value = -1030502.8

$ws.Range("H131").Value = 806.6
$ws.Range("I131").Value = 303.625
$ws.Range("J131").Value = 850.337
$ws.Range("K131").Value = 910.875
$ws.Range("L131").Value = 2551.011
$ws.Range("M131").Value = 4129.125
$ws.Range("N131").Value = -12631.011

$ws.Range("H132").Value = 1250
$ws.Range("I132").Value = 860
$ws.Range("J132").Value = 1373.1578
$ws.Range("K132").Value = 7740
$ws.Range("L132").Value = 12358.4202
$ws.Range("M132").Value = -5210
$ws.Range("N132").Value = -17418.4202

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 166668860
$ws.Range("J80").Value = 2051.5
$ws.Range("L80").Value = 2051.5
$ws.Range("N80").Value = -4047.5

$ws.Range("H83").Value = 166668860
$ws.Range("J83").Value = 2051.5
$ws.Range("L83").Value = 10257.5
$ws.Range("N83").Value = -20241.5

$ws.Range("H102").Value = 2343.6667
$ws.Range("I102").Value = 2343.6667
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 2343.6667
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -721.6667000000002
$ws.Range("N102").ClearContents()

$ws.Range("H120").Value = 34305.75
$ws.Range("J120").Value = 34305.75
$ws.Range("L120").Value = 34305.75
$ws.Range("N120").Value = -43981.75

$ws.Range("H122").Value = 1333
$ws.Range("I122").Value = 1000
$ws.Range("J122").Value = 1999
$ws.Range("K122").Value = 3000
$ws.Range("L122").Value = 5997
$ws.Range("M122").Value = -550
$ws.Range("N122").Value = -10897

$ws.Range("H126").Value = 3121.7334
$ws.Range("I126").Value = 2725.611
$ws.Range("J126").Value = 3715.9167
$ws.Range("K126").Value = 8176.833
$ws.Range("L126").Value = 11147.7501
$ws.Range("M126").Value = -5706.833
$ws.Range("N126").Value = -16087.7501

$ws.Range("H132").Value = 1478.0857
$ws.Range("I132").Value = 1028.2693
$ws.Range("J132").Value = 2777.5557
$ws.Range("K132").Value = 3084.8079
$ws.Range("L132").Value = 8332.667099999999
$ws.Range("M132").Value = -554.8078999999998
$ws.Range("N132").Value = -13392.6671

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 916.8823
$ws.Range("I22").Value = 866.6667
$ws.Range("J22").Value = 927.6429000000001
$ws.Range("K22").Value = 866.6667
$ws.Range("L22").Value = 927.6429000000001
$ws.Range("M22").Value = -571.6667
$ws.Range("N22").Value = -1517.6429

$ws.Range("H27").Value = 916.8823
$ws.Range("I27").Value = 866.6667
$ws.Range("J27").Value = 927.6429000000001
$ws.Range("K27").Value = 866.6667
$ws.Range("L27").Value = 927.6429000000001
$ws.Range("M27").Value = -759.6667
$ws.Range("N27").Value = -1141.6429

$ws.Range("H40").Value = 85449.914
$ws.Range("I40").Value = 500750
$ws.Range("J40").Value = 2389.9
$ws.Range("K40").Value = 500750
$ws.Range("L40").Value = 2389.9
$ws.Range("M40").Value = -500614
$ws.Range("N40").Value = -2661.9

$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 182737.55
$ws.Range("I81").Value = 200858
$ws.Range("J81").Value = 167637.17
$ws.Range("K81").Value = 401716
$ws.Range("L81").Value = 335274.34
$ws.Range("M81").Value = -400655
$ws.Range("N81").Value = -337396.34

$ws.Range("H84").Value = 182737.55
$ws.Range("I84").Value = 200858
$ws.Range("J84").Value = 167637.17
$ws.Range("K84").Value = 2008580
$ws.Range("L84").Value = 1676371.7
$ws.Range("M84").Value = -2003276
$ws.Range("N84").Value = -1686979.7

$ws.Range("H122").Value = 1777.7778
$ws.Range("I122").Value = 1000
$ws.Range("J122").Value = 4500
$ws.Range("K122").Value = 3000
$ws.Range("L122").Value = 13500
$ws.Range("M122").Value = -550
$ws.Range("N122").Value = -18400
